# Book1.xlsx / Sheet1 — "Resistance Rate" column (C2:C14) used to hold
# fractions (0.76, 0.88, ...) displayed with a percentage number format
# (76%, 88%, ...). The edit replaces that with plain numbers (76, 88, ...)
# using the default/General number format, and also moves the saved
# worksheet selection from D10 to F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Multiply each fraction by 100 so the stored value matches what used to be
# displayed (e.g. 0.76 -> 76), rounding away binary floating point noise
# such as 0.56000000000000005 -> 56.
for ($r = 2; $r -le 14; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = [Math]::Round($cell.Value2 * 100, 6)
}

# Drop the percentage number format from those cells so they fall back to
# the workbook default (General) formatting.
$ws.Range("C2:C14").ClearFormats()

# The saved file shows the active selection on F7 instead of D10.
$ws.Range("F7").Select()
